$d = $word.ActiveDocument

$d.Content.Find.Execute("with 9.5", $true, $false, $false, $false, $false,
                         $true, 1, $false, "with 10", 2)
